$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.442919373512268
$ws.Range("B1").Value = 3.316734552383423
$ws.Range("C1").Value = 4.297641277313232
$ws.Range("D1").Value = 1.994369268417358
$ws.Range("E1").Value = 1.157947421073914
